$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Fix B62 to reuse the clean shared string (removes orphaned variant)
$ws2.Cells.Item(62,2).Value = "gia_tri_con_lai"

# Populate column C with explanations, in the exact order the original author typed them
# (rows 1-54, then 56-59, then 55, then 60-68) so the shared-string table layout matches.
$ws2.Cells.Item(1,3).Value = "1. Đất, 2. Nhà, 3. Vật kiến trúc (Nhập 1,2,3,… theo danh mục nhóm tài sản)"
$ws2.Cells.Item(2,3).Value = "Nhập mã loại tài sản theo nhóm tài sản"
$ws2.Cells.Item(3,3).Value = "Mã tài sản"
$ws2.Cells.Item(4,3).Value = "Tên tài sản"
$ws2.Cells.Item(5,3).Value = "Lý do tăng: Đầu tư, Mua sắm, …"
$ws2.Cells.Item(6,3).Value = "Số lượng"
$ws2.Cells.Item(7,3).Value = "Đơn vị tính"
$ws2.Cells.Item(8,3).Value = "Mã bộ phận sử dụng"
$ws2.Cells.Item(9,3).Value = "Mã tỉnh thành phố theo danh mục"
$ws2.Cells.Item(10,3).Value = "Mã thuyện, quận theo danh mục"
$ws2.Cells.Item(11,3).Value = "Mã xã theo danh mục"
$ws2.Cells.Item(12,3).Value = "Địa chỉ"
$ws2.Cells.Item(13,3).Value = "Số tầng"
$ws2.Cells.Item(14,3).Value = "chiều dài"
$ws2.Cells.Item(15,3).Value = "diện tích xây dựng"
$ws2.Cells.Item(16,3).Value = "thể tích"
$ws2.Cells.Item(17,3).Value = "năm xây dựng"
$ws2.Cells.Item(18,3).Value = "nước sản xuất"
$ws2.Cells.Item(19,3).Value = "biển kiểm soát"
$ws2.Cells.Item(20,3).Value = "nhãn hiệu tài sản"
$ws2.Cells.Item(21,3).Value = "Model"
$ws2.Cells.Item(22,3).Value = "số seri"
$ws2.Cells.Item(23,3).Value = "số máy"
$ws2.Cells.Item(24,3).Value = "tải trọng"
$ws2.Cells.Item(25,3).Value = "số chỗ ngồi"
$ws2.Cells.Item(26,3).Value = "số cầu"
$ws2.Cells.Item(27,3).Value = "công suất "
$ws2.Cells.Item(28,3).Value = "dung tích xe"
$ws2.Cells.Item(29,3).Value = "giấy công nhận đăng ký số"
$ws2.Cells.Item(30,3).Value = "ngày đăng ký"
$ws2.Cells.Item(31,3).Value = "cơ quan cấp đăng ký"
$ws2.Cells.Item(32,3).Value = "nguồn gốc xe"
$ws2.Cells.Item(33,3).Value = "màu sơn"
$ws2.Cells.Item(34,3).Value = "người sử dụng"
$ws2.Cells.Item(35,3).Value = "hình thức bố trí sử dụng"
$ws2.Cells.Item(36,3).Value = "chức danh sử dụng"
$ws2.Cells.Item(37,3).Value = "quyết định trang cấp"
$ws2.Cells.Item(38,3).Value = "ngày quyết định trang cấp"
$ws2.Cells.Item(39,3).Value = "dự án"
$ws2.Cells.Item(40,3).Value = "1. Đất, 2. Nhà, 3. Xe ô tô, 4. Tài sản trên 500tr, 5. tài sản dưới 500tr"
$ws2.Cells.Item(41,3).Value = "thông số kỹ thuật"
$ws2.Cells.Item(42,3).Value = "quản lý nhà nước (hoạt động sự nghiệp)"
$ws2.Cells.Item(43,3).Value = "HĐSN - Không KD"
$ws2.Cells.Item(44,3).Value = "HĐSN - KD"
$ws2.Cells.Item(45,3).Value = "HĐSN - LDLK"
$ws2.Cells.Item(46,3).Value = "HĐSN - Cho thuê"
$ws2.Cells.Item(47,3).Value = "sử dụng khác"
$ws2.Cells.Item(48,3).Value = "trạng thái: 0. chưa ghi tăng, 1. đang sử dụng"
$ws2.Cells.Item(49,3).Value = "tổng diện tích"
$ws2.Cells.Item(50,3).Value = "giá trị đất"
$ws2.Cells.Item(51,3).Value = "ngày mua"
$ws2.Cells.Item(52,3).Value = "ngày bắt đầu sử dụng"
$ws2.Cells.Item(53,3).Value = "ngày ghi tăng"
$ws2.Cells.Item(54,3).Value = "năm theo dõi"
$ws2.Cells.Item(56,3).Value = "số năm sử dụng"
$ws2.Cells.Item(57,3).Value = "tỷ lệ hao mòn"
$ws2.Cells.Item(58,3).Value = "HM/KH năm"
$ws2.Cells.Item(59,3).Value = "số năm sử dụng còn lại"
$ws2.Cells.Item(55,3).Value = "ngày bắt đầu tính hao mòn yyyy-MM-dd, vd 2024-12-30"
$ws2.Cells.Item(60,3).Value = "ngày kết thúc hao mòn yyyy-MM-dd v,  2024-12-30"
$ws2.Cells.Item(61,3).Value = "hao mòn lũy kế"
$ws2.Cells.Item(62,3).Value = "giá trị còn lại"
$ws2.Cells.Item(63,3).Value = "mục đích sử dụng: Đất hoạt động sự nghiệp …"
$ws2.Cells.Item(64,3).Value = "diện tích để ở"
$ws2.Cells.Item(65,3).Value = "dt bỏ trống"
$ws2.Cells.Item(66,3).Value = "dt bị lấn chiếm"
$ws2.Cells.Item(67,3).Value = "dt sử dụng hỗn hợp"
$ws2.Cells.Item(68,3).Value = "mặc định: vpddt"

# Size column C to fit the explanatory text (closest achievable to saved width 86.5703125)
$ws2.Columns.Item(3).ColumnWidth = 85.67

# Make GiaiThich the active sheet/tab and mirror the saved selection/scroll state
$ws2.Activate()
$ws2.Range("C70").Select() | Out-Null
